$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 110 (this pushes the old row 110 and
# everything below it down by one row, to 111, 112, ... 212).
$ws.Rows.Item(110).Insert()

# The new row 110 shares every "descriptive" column (Mercado ID, Mercado,
# Region, Codreg, Categoria ID, Categoria, Variedad, Calidad, Unidad de
# comercializacion, Origen, Kg o Unidades, Clasificacion) with the record
# that is now directly below it (row 111, the old row 110) - copy those
# across first, then overwrite the few cells that actually hold new data.
for ($col = 1; $col -le 18; $col++) {
    $ws.Cells.Item(110, $col).Value = $ws.Cells.Item(111, $col).Value2
}

$ws.Cells.Item(110, 4).Value  = 44589  # D110 - Fecha
$ws.Cells.Item(110, 10).Value = 240    # J110 - Volumen
$ws.Cells.Item(110, 11).Value = 20000  # K110 - Precio minimo
$ws.Cells.Item(110, 12).Value = 21000  # L110 - Precio maximo
$ws.Cells.Item(110, 13).Value = 20500  # M110 - Precio promedio ponderado
$ws.Cells.Item(110, 16).Value = 2050   # P110 - Precio $/Kg
